$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that are plain text (including European-style thousand-separated
# numbers and percentage strings) - assigned directly since Excel will not
# mis-detect them as numbers.
$textUpdates = @{
    'D2' = '27.938.70'
    'E2' = '  +1.44%  '
    'D3' = '1.642.88'
    'E3' = '  +1.18%  '
    'E4' = '  -0.03%  '
    'E5' = '  +0.87%  '
    'E6' = '  -0.06%  '
    'E7' = '  -0.05%  '
    'E8' = '  +2.82%  '
    'E9' = '  +0.64%  '
    'E10' = '  +0.79%  '
    'E11' = '  -0.58%  '
    'D12' = '1.876.00'
    'E12' = '  +1.20%  '
    'D13' = '1.646.06'
    'E13' = '  +1.31%  '
    'E14' = '  +4.81%  '
    'E15' = '  +0.96%  '
    'E16' = '  +1.19%  '
    'D17' = '27.924.38'
    'E17' = '  +1.53%  '
    'E18' = '  +0.19%  '
    'E19' = '  +0.94%  '
    'E20' = '  +1.33%  '
    'E21' = '  -0.05%  '
    'E22' = '  +6.96%  '
    'E23' = '  +1.54%  '
    'E24' = '  -0.34%  '
    'E25' = '  +2.57%  '
    'E26' = '  +0.87%  '
    'E27' = '  +0.72%  '
    'E28' = '  +1.27%  '
    'E29' = '  -0.06%  '
    'E30' = '  +1.37%  '
    'E31' = '  +0.46%  '
    'E32' = '  +2.14%  '
    'D33' = '1.422.97'
    'E33' = '  -3.15%  '
    'E34' = '  +2.04%  '
    'E35' = '  +2.01%  '
    'E36' = '  +0.12%  '
    'E37' = '  +2.18%  '
    'E38' = '  -0.46%  '
    'E39' = '  +1.01%  '
    'E40' = '  +0.61%  '
    'E41' = '  +2.32%  '
    'E42' = '  -0.06%  '
    'E43' = '  +0.28%  '
    'E44' = '  +0.39%  '
    'E45' = '  +3.15%  '
    'E46' = '  +3.09%  '
    'E47' = '  -0.01%  '
    'D48' = '1.784.77'
    'E48' = '  +1.19%  '
    'E49' = '  +1.92%  '
    'E50' = '  +1.10%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E51' = '  +0.73%  '
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}

# Cell values that look like plain numbers (e.g. "213.54") need to be forced
# to text format first, otherwise Excel auto-converts them into numeric cells
# (losing the original text representation / trailing zero formatting).
$numericLookingUpdates = @{
    'D5' = '213.54'
    'D8' = '23.88'
    'D11' = '0.0876'
    'D14' = '0.576'
    'D16' = '65.91'
    'D18' = '230.75'
    'D22' = '11.08'
    'D24' = '2.07'
    'D25' = '152.66'
    'D28' = '15.74'
    'D32' = '3.34'
    'D34' = '3.12'
    'D37' = '0.892'
    'D38' = '0.928'
    'D40' = '0.559'
    'D43' = '67.31'
    'D44' = '2.47'
    'D45' = '5.45'
    'D49' = '88.93'
    'D51' = '7.71'
}

foreach ($cellRef in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$cellRef]
}
